$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 1.249563
$ws.Cells.Item(2, 8).Value = 3.748689
$ws.Cells.Item(2, 9).Value = 0.006879661005457268
$ws.Cells.Item(2, 10).Value = 0.006920555383629864
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.506715
$ws.Cells.Item(2, 14).Value = 1.520145
$ws.Cells.Item(2, 15).Value = 0.003122343715987576
$ws.Cells.Item(2, 16).Value = 0.003132472094339857
$ws.Cells.Item(2, 17).Value = 0.633172315545
$ws.Cells.Item(2, 18).Value = 5.698550839905
$ws.Cells.Item(2, 19).Value = [double]"2.148066630851427E-05"
$ws.Cells.Item(2, 20).Value = [double]"2.167844661655401E-05"

# Row 3
$ws.Cells.Item(3, 7).Value = 1.249563
$ws.Cells.Item(3, 8).Value = 3.748689
$ws.Cells.Item(3, 9).Value = 0.006879661005457268
$ws.Cells.Item(3, 10).Value = 0.006920555383629864
$ws.Cells.Item(3, 13).Value = 88.13219433333332
$ws.Cells.Item(3, 14).Value = 264.396583
$ws.Cells.Item(3, 15).Value = 0.5430646480820168
$ws.Cells.Item(3, 16).Value = 0.5448262620252092
$ws.Cells.Item(3, 17).Value = 110.126729147743
$ws.Cells.Item(3, 18).Value = 991.1405623296868
$ws.Cells.Item(3, 19).Value = 0.003736100682852225
$ws.Cells.Item(3, 20).Value = 0.003770500320801497

# Row 4
$ws.Cells.Item(4, 7).Value = 1.249563
$ws.Cells.Item(4, 8).Value = 3.748689
$ws.Cells.Item(4, 9).Value = 0.006879661005457268
$ws.Cells.Item(4, 10).Value = 0.006920555383629864
$ws.Cells.Item(4, 13).Value = 1.5741895
$ws.Cells.Item(4, 14).Value = 3.148379
$ws.Cells.Item(4, 15).Value = 0.009700049718478087
$ws.Cells.Item(4, 16).Value = 0.006487676741301404
$ws.Cells.Item(4, 17).Value = 1.9670489541885
$ws.Cells.Item(4, 18).Value = 11.802293725131
$ws.Cells.Item(4, 19).Value = [double]"6.673305379921045E-05"
$ws.Cells.Item(4, 20).Value = [double]"4.489832619926369E-05"

# Row 5
$ws.Cells.Item(5, 7).Value = 1.249563
$ws.Cells.Item(5, 8).Value = 3.748689
$ws.Cells.Item(5, 9).Value = 0.006879661005457268
$ws.Cells.Item(5, 10).Value = 0.006920555383629864
$ws.Cells.Item(5, 13).Value = 72.07364666666666
$ws.Cells.Item(5, 14).Value = 216.22094
$ws.Cells.Item(5, 15).Value = 0.4441129584835175
$ws.Cells.Item(5, 16).Value = 0.4455535891391496
$ws.Cells.Item(5, 17).Value = 90.06056214973999
$ws.Cells.Item(5, 18).Value = 810.5450593476598
$ws.Cells.Item(5, 19).Value = 0.003055346602497318
$ws.Cells.Item(5, 20).Value = 0.00308347829001255

# Row 6
$ws.Cells.Item(6, 8).Value = 531.128907
$ws.Cells.Item(6, 9).Value = 0.9747372562405255
$ws.Cells.Item(6, 10).Value = 0.9805313315509224
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.506715
$ws.Cells.Item(6, 14).Value = 1.520145
$ws.Cells.Item(6, 15).Value = 0.003122343715987576
$ws.Cells.Item(6, 16).Value = 0.003132472094339857
$ws.Cells.Item(6, 17).Value = 89.71032803683501
$ws.Cells.Item(6, 18).Value = 807.3929523315151
$ws.Cells.Item(6, 19).Value = 0.003043464746761577
$ws.Cells.Item(6, 20).Value = 0.003071487033709166

# Row 7
$ws.Cells.Item(7, 8).Value = 531.128907
$ws.Cells.Item(7, 9).Value = 0.9747372562405255
$ws.Cells.Item(7, 10).Value = 0.9805313315509224
$ws.Cells.Item(7, 13).Value = 88.13219433333332
$ws.Cells.Item(7, 14).Value = 264.396583
$ws.Cells.Item(7, 15).Value = 0.5430646480820168
$ws.Cells.Item(7, 16).Value = 0.5448262620252092
$ws.Cells.Item(7, 17).Value = 15603.18534925831
$ws.Cells.Item(7, 18).Value = 140428.6681433248
$ws.Cells.Item(7, 19).Value = 0.5293453450326917
$ws.Cells.Item(7, 20).Value = 0.5342192201674901

# Row 8
$ws.Cells.Item(8, 8).Value = 531.128907
$ws.Cells.Item(8, 9).Value = 0.9747372562405255
$ws.Cells.Item(8, 10).Value = 0.9805313315509224
$ws.Cells.Item(8, 13).Value = 1.5741895
$ws.Cells.Item(8, 14).Value = 3.148379
$ws.Cells.Item(8, 15).Value = 0.009700049718478087
$ws.Cells.Item(8, 16).Value = 0.006487676741301404
$ws.Cells.Item(8, 17).Value = 278.6991828486255
$ws.Cells.Item(8, 18).Value = 1672.195097091753
$ws.Cells.Item(8, 19).Value = 0.009454999847986012
$ws.Cells.Item(8, 20).Value = 0.006361370313820215

# Row 9
$ws.Cells.Item(9, 8).Value = 531.128907
$ws.Cells.Item(9, 9).Value = 0.9747372562405255
$ws.Cells.Item(9, 10).Value = 0.9805313315509224
$ws.Cells.Item(9, 13).Value = 72.07364666666666
$ws.Cells.Item(9, 14).Value = 216.22094
$ws.Cells.Item(9, 15).Value = 0.4441129584835175
$ws.Cells.Item(9, 16).Value = 0.4455535891391496
$ws.Cells.Item(9, 17).Value = 12760.13239252362
$ws.Cells.Item(9, 18).Value = 114841.1915327126
$ws.Cells.Item(9, 19).Value = 0.4328934466130862
$ws.Cells.Item(9, 20).Value = 0.4368792540359029

# Row 10
$ws.Cells.Item(10, 7).Value = 0.07049699999999999
$ws.Cells.Item(10, 8).Value = 0.211491
$ws.Cells.Item(10, 9).Value = 0.000388132060489724
$ws.Cells.Item(10, 10).Value = 0.0003904392118522671
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.506715
$ws.Cells.Item(10, 14).Value = 1.520145
$ws.Cells.Item(10, 15).Value = 0.003122343715987576
$ws.Cells.Item(10, 16).Value = 0.003132472094339857
$ws.Cells.Item(10, 17).Value = 0.035721887355
$ws.Cells.Item(10, 18).Value = 0.321496986195
$ws.Cells.Item(10, 19).Value = [double]"1.2118817000434E-06"
$ws.Cells.Item(10, 20).Value = [double]"1.223039935663274E-06"

# Row 11
$ws.Cells.Item(11, 7).Value = 0.07049699999999999
$ws.Cells.Item(11, 8).Value = 0.211491
$ws.Cells.Item(11, 9).Value = 0.000388132060489724
$ws.Cells.Item(11, 10).Value = 0.0003904392118522671
$ws.Cells.Item(11, 13).Value = 88.13219433333332
$ws.Cells.Item(11, 14).Value = 264.396583
$ws.Cells.Item(11, 15).Value = 0.5430646480820168
$ws.Cells.Item(11, 16).Value = 0.5448262620252092
$ws.Cells.Item(11, 17).Value = 6.213055303916998
$ws.Cells.Item(11, 18).Value = 55.91749773525299
$ws.Cells.Item(11, 19).Value = 0.0002107808008392
$ws.Cells.Item(11, 20).Value = 0.0002127215363415395

# Row 12
$ws.Cells.Item(12, 7).Value = 0.07049699999999999
$ws.Cells.Item(12, 8).Value = 0.211491
$ws.Cells.Item(12, 9).Value = 0.000388132060489724
$ws.Cells.Item(12, 10).Value = 0.0003904392118522671
$ws.Cells.Item(12, 13).Value = 1.5741895
$ws.Cells.Item(12, 14).Value = 3.148379
$ws.Cells.Item(12, 15).Value = 0.009700049718478087
$ws.Cells.Item(12, 16).Value = 0.006487676741301404
$ws.Cells.Item(12, 17).Value = 0.1109756371815
$ws.Cells.Item(12, 18).Value = 0.6658538230889999
$ws.Cells.Item(12, 19).Value = [double]"3.764900284085667E-06"
$ws.Cells.Item(12, 20).Value = [double]"2.533043393626005E-06"

# Row 13
$ws.Cells.Item(13, 7).Value = 0.07049699999999999
$ws.Cells.Item(13, 8).Value = 0.211491
$ws.Cells.Item(13, 9).Value = 0.000388132060489724
$ws.Cells.Item(13, 10).Value = 0.0003904392118522671
$ws.Cells.Item(13, 13).Value = 72.07364666666666
$ws.Cells.Item(13, 14).Value = 216.22094
$ws.Cells.Item(13, 15).Value = 0.4441129584835175
$ws.Cells.Item(13, 16).Value = 0.4455535891391496
$ws.Cells.Item(13, 17).Value = 5.080975869059999
$ws.Cells.Item(13, 18).Value = 45.72878282153999
$ws.Cells.Item(13, 19).Value = 0.0001723744776663949
$ws.Cells.Item(13, 20).Value = 0.0001739615921814384

# Row 14
$ws.Cells.Item(14, 7).Value = 3.2198455
$ws.Cells.Item(14, 8).Value = 6.439691
$ws.Cells.Item(14, 9).Value = 0.01772735390688349
$ws.Cells.Item(14, 10).Value = 0.01188848640657115
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.506715
$ws.Cells.Item(14, 14).Value = 1.520145
$ws.Cells.Item(14, 15).Value = 0.003122343715987576
$ws.Cells.Item(14, 16).Value = 0.003132472094339857
$ws.Cells.Item(14, 17).Value = 1.6315440125325
$ws.Cells.Item(14, 18).Value = 9.789264075195
$ws.Cells.Item(14, 19).Value = [double]"5.535089207224548E-05"
$ws.Cells.Item(14, 20).Value = [double]"3.724035191252284E-05"

# Row 15
$ws.Cells.Item(15, 7).Value = 3.2198455
$ws.Cells.Item(15, 8).Value = 6.439691
$ws.Cells.Item(15, 9).Value = 0.01772735390688349
$ws.Cells.Item(15, 10).Value = 0.01188848640657115
$ws.Cells.Item(15, 13).Value = 88.13219433333332
$ws.Cells.Item(15, 14).Value = 264.396583
$ws.Cells.Item(15, 15).Value = 0.5430646480820168
$ws.Cells.Item(15, 16).Value = 0.5448262620252092
$ws.Cells.Item(15, 17).Value = 283.7720493293088
$ws.Cells.Item(15, 18).Value = 1702.632295975853
$ws.Cells.Item(15, 19).Value = 0.00962709921086705
$ws.Cells.Item(15, 20).Value = 0.00647715961002967

# Row 16
$ws.Cells.Item(16, 7).Value = 3.2198455
$ws.Cells.Item(16, 8).Value = 6.439691
$ws.Cells.Item(16, 9).Value = 0.01772735390688349
$ws.Cells.Item(16, 10).Value = 0.01188848640657115
$ws.Cells.Item(16, 13).Value = 1.5741895
$ws.Cells.Item(16, 14).Value = 3.148379
$ws.Cells.Item(16, 15).Value = 0.009700049718478087
$ws.Cells.Item(16, 16).Value = 0.006487676741301404
$ws.Cells.Item(16, 17).Value = 5.068646977722249
$ws.Cells.Item(16, 18).Value = 20.274587910889
$ws.Cells.Item(16, 19).Value = 0.0001719562142738266
$ws.Cells.Item(16, 20).Value = [double]"7.712865674918953E-05"

# Row 17
$ws.Cells.Item(17, 7).Value = 3.2198455
$ws.Cells.Item(17, 8).Value = 6.439691
$ws.Cells.Item(17, 9).Value = 0.01772735390688349
$ws.Cells.Item(17, 10).Value = 0.01188848640657115
$ws.Cells.Item(17, 13).Value = 72.07364666666666
$ws.Cells.Item(17, 14).Value = 216.22094
$ws.Cells.Item(17, 15).Value = 0.4441129584835175
$ws.Cells.Item(17, 16).Value = 0.4455535891391496
$ws.Cells.Item(17, 17).Value = 232.0660068882567
$ws.Cells.Item(17, 18).Value = 1392.39604132954
$ws.Cells.Item(17, 19).Value = 0.00787294758967037
$ws.Cells.Item(17, 20).Value = 0.005296957787879765

# Row 18
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 0.6666666666666666
$ws.Cells.Item(18, 7).Value = 0.048604
$ws.Cells.Item(18, 8).Value = 0.145812
$ws.Cells.Item(18, 9).Value = 0.0002675967866440068
$ws.Cells.Item(18, 10).Value = 0.0002691874470242364
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.506715
$ws.Cells.Item(18, 14).Value = 1.520145
$ws.Cells.Item(18, 15).Value = 0.003122343715987576
$ws.Cells.Item(18, 16).Value = 0.003132472094339857
$ws.Cells.Item(18, 17).Value = 0.02462837586
$ws.Cells.Item(18, 18).Value = 0.22165538274
$ws.Cells.Item(18, 19).Value = [double]"8.355291451963826E-07"
$ws.Cells.Item(18, 20).Value = [double]"8.432221659500091E-07"

# Row 19
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 0.6666666666666666
$ws.Cells.Item(19, 7).Value = 0.048604
$ws.Cells.Item(19, 8).Value = 0.145812
$ws.Cells.Item(19, 9).Value = 0.0002675967866440068
$ws.Cells.Item(19, 10).Value = 0.0002691874470242364
$ws.Cells.Item(19, 13).Value = 88.13219433333332
$ws.Cells.Item(19, 14).Value = 264.396583
$ws.Cells.Item(19, 15).Value = 0.5430646480820168
$ws.Cells.Item(19, 16).Value = 0.5448262620252092
$ws.Cells.Item(19, 17).Value = 4.283577173377332
$ws.Cells.Item(19, 18).Value = 38.552194560396
$ws.Cells.Item(19, 19).Value = 0.0001453223547667061
$ws.Cells.Item(19, 20).Value = 0.0001466603905463238

# Row 20
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 0.6666666666666666
$ws.Cells.Item(20, 7).Value = 0.048604
$ws.Cells.Item(20, 8).Value = 0.145812
$ws.Cells.Item(20, 9).Value = 0.0002675967866440068
$ws.Cells.Item(20, 10).Value = 0.0002691874470242364
$ws.Cells.Item(20, 13).Value = 1.5741895
$ws.Cells.Item(20, 14).Value = 3.148379
$ws.Cells.Item(20, 15).Value = 0.009700049718478087
$ws.Cells.Item(20, 16).Value = 0.006487676741301404
$ws.Cells.Item(20, 17).Value = 0.076511906458
$ws.Cells.Item(20, 18).Value = 0.4590714387479999
$ws.Cells.Item(20, 19).Value = [double]"2.595702134951838E-06"
$ws.Cells.Item(20, 20).Value = [double]"1.746401139109443E-06"

# Row 21
$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 0.6666666666666666
$ws.Cells.Item(21, 7).Value = 0.048604
$ws.Cells.Item(21, 8).Value = 0.145812
$ws.Cells.Item(21, 9).Value = 0.0002675967866440068
$ws.Cells.Item(21, 10).Value = 0.0002691874470242364
$ws.Cells.Item(21, 13).Value = 72.07364666666666
$ws.Cells.Item(21, 14).Value = 216.22094
$ws.Cells.Item(21, 15).Value = 0.4441129584835175
$ws.Cells.Item(21, 16).Value = 0.4455535891391496
$ws.Cells.Item(21, 17).Value = 3.503067522586667
$ws.Cells.Item(21, 18).Value = 31.52760770328
$ws.Cells.Item(21, 19).Value = 0.0001188432005971525
$ws.Cells.Item(21, 20).Value = 0.0001199374331728532
